# Daily attendance processing - 2025-11-30 03:22:00
#
# Normalizes the "Recorded By" column (G): for rows recorded by one of the
# known automated accounts (dnasr281@gmail.com or backup@backdoor.com)
# together with "System", move the "System" token so that it appears right
# after the initial lowercase "system" marker (if present), or otherwise
# move it to the very front of the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = [string]$cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    if (-not ($text.Contains("backup@backdoor.com") -or $text.Contains("dnasr281@gmail.com"))) {
        continue
    }

    $rawParts = $text.Split(",")
    $trimmedParts = @()
    foreach ($p in $rawParts) {
        $trimmedParts += $p.Trim()
    }

    $hasSystem = $false
    foreach ($p in $trimmedParts) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        }
    }

    if (-not $hasSystem) {
        continue
    }

    $filteredParts = @()
    foreach ($p in $trimmedParts) {
        if (-not $p.Equals("System")) {
            $filteredParts += $p
        }
    }

    if ($filteredParts.Length -gt 0 -and $filteredParts[0].Equals("system")) {
        if ($filteredParts.Length -gt 1) {
            $newParts = @($filteredParts[0], "System") + $filteredParts[1..($filteredParts.Length - 1)]
        } else {
            $newParts = @($filteredParts[0], "System")
        }
    } else {
        $newParts = @("System") + $filteredParts
    }

    $newText = $newParts -join ", "

    if (-not $newText.Equals($text)) {
        $cell.Value = $newText
    }
}
